$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (they are decimal-looking strings
# that must stay as literal text, matching the source data feed format).
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.508.10"
$ws.Range("D3").Value = "2.109.06"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "334.78"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "0.5236"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").Value = "0.4537"
$ws.Range("E8").Value = "  +3.59%  "
$ws.Range("D9").Value = "53.42"
$ws.Range("E9").Value = "  +13.65%  "
$ws.Range("D10").Value = "0.08996"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "1.163"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "24.47"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").Value = "2.103.63"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "6.779"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "7.823"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "0.06619"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "19.33"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "6.310"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "30.565.72"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "12.38"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "2.358"
$ws.Range("E25").Value = "  +3.87%  "
$ws.Range("D26").Value = "2.354.84"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").Value = "22.38"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "2.579"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "163.31"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "132.84"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "1.200"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").Value = "0.1075"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "1.662"
$ws.Range("E33").Value = "  +6.42%  "
$ws.Range("D34").Value = "6.167"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "3.947"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "10.64"
$ws.Range("E36").Value = "  +12.19%  "
$ws.Range("D37").Value = "0.02578"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "0.06810"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "5.546"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "12.77"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "0.2293"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "0.6926"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "1.258"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "2.391"
$ws.Range("E44").Value = "  +7.28%  "
$ws.Range("D46").Value = "0.6415"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "14.04"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "3.664"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "1.216"
$ws.Range("E50").Value = "  +4.76%  "
$ws.Range("D51").Value = "83.42"
$ws.Range("E51").Value = "  +0.41%  "
